$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288, shifting existing rows 288..366 down to 289..367
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new record
$ws.Cells.Item(288, 1).Value = 3
$ws.Cells.Item(288, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(288, 3).Value = "Coquimbo"
$ws.Cells.Item(288, 4).Value = 44551
$ws.Cells.Item(288, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(288, 5).Value = 5
$ws.Cells.Item(288, 6).Value = 100112037
$ws.Cells.Item(288, 7).Value = "Cebollín"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 310
$ws.Cells.Item(288, 11).Value = 3000
$ws.Cells.Item(288, 12).Value = 3500
$ws.Cells.Item(288, 13).Value = 3258
$ws.Cells.Item(288, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(288, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(288, 16).Value = 90
$ws.Cells.Item(288, 17).Value = 36
$ws.Cells.Item(288, 18).Value = "Hortaliza"
